$wb = $excel.ActiveWorkbook

# Add the new "Logging" sheet (inserted before the active sheet)
$logging = $wb.Worksheets.Add()
$logging.Name = "Logging"
$logging.Range("B1").Value = "carryover"
$logging.Range("A2").Value = "row"
$logging.Range("B2").Value = 36
$logging.Range("A3").Value = "column"
$logging.Range("B3").Value = 10
$logging.Activate()
$logging.Range("B2").Select() | Out-Null

# Rename the original sheet to "Timesheet"
$ts = $wb.Worksheets.Item(2)
$ts.Name = "Timesheet"

# Make Timesheet the active sheet/tab, with the new selection
$ts.Activate()
$ts.Range("J36").Select() | Out-Null
